$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 533
$ws.Range("L3").Value = 538
$ws.Range("B4").Value = 1704
$ws.Range("K4").Value = 1733
$ws.Range("L4").Value = 140
$ws.Range("L5").Value = 45
$ws.Range("J6").Value = 11052
$ws.Range("L6").Value = 599
$ws.Range("B7").Value = 23337
$ws.Range("J7").Value = 29320
$ws.Range("K7").Value = 27522
$ws.Range("L7").Value = 1855

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L6").Value = 9
$ws.Range("L7").Value = 27

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 28
$ws.Range("L3").Value = 31
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 101

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 11
$ws.Range("L7").Value = 37

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 12
$ws.Range("L3").Value = 30
$ws.Range("J6").Value = 473
$ws.Range("J7").Value = 1312
$ws.Range("L7").Value = 77

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L6").Value = 17
$ws.Range("L7").Value = 59

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L3").Value = 14
$ws.Range("L7").Value = 32

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("L6").Value = 5
$ws.Range("L7").Value = 12

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L7").Value = 65
$ws.Range("L8").Value = 101
$ws.Range("L19").Value = 63
$ws.Range("L20").Value = 50
$ws.Range("L21").Value = 5
$ws.Range("L22").Value = 6
$ws.Range("L23").Value = 17
$ws.Range("L26").Value = 3
$ws.Range("L27").Value = 15
$ws.Range("L29").Value = 101
$ws.Range("L30").Value = 12
$ws.Range("J33").Value = 1312
$ws.Range("L33").Value = 77
$ws.Range("L37").Value = 59
$ws.Range("L42").Value = 67
$ws.Range("L46").Value = 7
$ws.Range("L48").Value = 32
$ws.Range("L51").Value = 26
$ws.Range("L52").Value = 36
$ws.Range("L53").Value = 27
$ws.Range("L54").Value = 33
$ws.Range("B63").Value = 408
$ws.Range("L63").Value = 12
$ws.Range("L67").Value = 55
$ws.Range("L71").Value = 6
$ws.Range("L76").Value = 28
$ws.Range("L78").Value = 19
$ws.Range("L79").Value = 51
$ws.Range("L83").Value = 37
$ws.Range("L85").Value = 94
$ws.Range("K88").Value = 289
$ws.Range("L88").Value = 24
$ws.Range("K89").Value = 411
$ws.Range("L89").Value = 20
$ws.Range("L93").Value = 11
$ws.Range("L94").Value = 26
$ws.Range("L97").Value = 25
$ws.Range("L98").Value = 16
$ws.Range("L99").Value = 32
$ws.Range("B101").Value = 23337
$ws.Range("J101").Value = 29320
$ws.Range("K101").Value = 27522
$ws.Range("L101").Value = 1855

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 15
$ws.Range("L3").Value = 14
$ws.Range("L7").Value = 55

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L2").Value = 10
$ws.Range("L7").Value = 33

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 32
$ws.Range("L3").Value = 30
$ws.Range("L7").Value = 101

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L2").Value = 5
$ws.Range("L3").Value = 6
$ws.Range("L7").Value = 32

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 24
$ws.Range("L3").Value = 14
$ws.Range("L6").Value = 22
$ws.Range("L7").Value = 63

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L6").Value = 17
$ws.Range("L7").Value = 28

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 15
$ws.Range("L5").Value = 2
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 67

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L2").Value = 5
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 19

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("L4").Value = 1
$ws.Range("L7").Value = 7

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L3").Value = 8
$ws.Range("L7").Value = 17

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("L6").Value = 4
$ws.Range("L7").Value = 5

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 20
$ws.Range("L7").Value = 51

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 15
$ws.Range("L3").Value = 14
$ws.Range("L7").Value = 50

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L2").Value = 5
$ws.Range("L7").Value = 11

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L6").Value = 24
$ws.Range("L7").Value = 65

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L2").Value = 7
$ws.Range("L7").Value = 26

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 16

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("L6").Value = 3
$ws.Range("L7").Value = 3

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L6").Value = 22
$ws.Range("L7").Value = 25

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K4").Value = 8
$ws.Range("L4").Value = 2
$ws.Range("L6").Value = 10
$ws.Range("K7").Value = 289
$ws.Range("L7").Value = 24

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L2").Value = 9
$ws.Range("K4").Value = 49
$ws.Range("K7").Value = 411
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L6").Value = 5
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L2").Value = 8
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 26

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 23
$ws.Range("L3").Value = 42
$ws.Range("L7").Value = 94

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("L3").Value = 2
$ws.Range("L6").Value = 2
$ws.Range("L7").Value = 6

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("L2").Value = 2
$ws.Range("L7").Value = 6

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 14
$ws.Range("L3").Value = 10
$ws.Range("L7").Value = 36
